# hysplits.xlsx direction-code fix
# Commit: "fixed hysplits.xlsx to INP(T) data, fixed Niemand param to cover
#          all data in contour3, fixed graph in hysplits.py to give viewable
#          legend size"
#
# Columns C (Scan) and D (Irl) hold compass-direction letters (N/S/E/W).
# A batch of rows had the wrong letters recorded; this script corrects
# them to match the real INP(T) data, adds a couple of missing D-column
# entries, flags the corrected Niemand-param block (rows 54-57) in red,
# moves the active selection, and sets the print page to A4 portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-5: Scan (C) was "N", should be "S" ---
$ws.Range("C2:C5").Value = "S"

# --- Rows 10-13: Scan (C) was "W", should be "S"; Irl (D) gains "W" ---
$ws.Range("C10:C13").Value = "S"
$ws.Range("D10").Value = "W"
$ws.Range("D11").Value = "W"
$ws.Range("D12").Value = "W"
$ws.Range("D13").Value = "W"

# --- Rows 26-29: Scan/Irl were swapped (E/N -> N/E) ---
$ws.Range("C26:C29").Value = "N"
$ws.Range("D26:D29").Value = "E"

# --- Rows 54-57: Niemand param block, corrected + highlighted red ---
$ws.Range("C54").Value = "S"
$ws.Range("D54").Value = "W"
$ws.Range("C55").Value = "S"
$ws.Range("D55").Value = "W"
$ws.Range("C56").Value = "S"
$ws.Range("D56").Value = "W"
$ws.Range("C57").Value = "S"

$red = $ws.Range("C54:C57")
$red.Font.Color = 255
$red.HorizontalAlignment = -4108

# --- Rows 58-61: Irl (D) was missing, add "N" ---
$ws.Range("D58").Value = "N"
$ws.Range("D59").Value = "N"
$ws.Range("D60").Value = "N"
$ws.Range("D61").Value = "N"

# --- Rows 62-65: Scan (C) was "E", should be "N" ---
$ws.Range("C62:C65").Value = "N"

# --- Rows 78-81: Irl (D) was missing, add "S" ---
$ws.Range("D78").Value = "S"
$ws.Range("D79").Value = "S"
$ws.Range("D80").Value = "S"
$ws.Range("D81").Value = "S"

# --- Rows 85-91: Scan (C) was "E"/"W", should be "S" ---
$ws.Range("C85:C89").Value = "S"
$ws.Range("C90").Value = "S"
$ws.Range("C91").Value = "S"

# --- Move active selection to C81 (last cell touched while editing) ---
$ws.Range("C81").Select()

# --- Print setup: A4, portrait ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
